$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.53%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'9.44%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.267"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.20%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07464"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'7.20%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.841"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.63%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.794"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'6.73%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.522"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'8.46%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'2.14%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1682"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.82%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08008"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'5.79%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07969"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.27%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03006"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.56%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09898"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'9.83%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001500"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-5.40%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.04601"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.54%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006484"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.86%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.457"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.90%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.230"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.00%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E21").Value = "'0.05%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.494"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'10.99%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'1.32%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.54%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004444"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'7.23%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001399"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'19.62%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001691"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'1.03%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01725"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2,541.26%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04489"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.55%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007152"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.04%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1348"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'8.04%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002139"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.35%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01277"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'9.20%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006170"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.64%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.7093"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-63.24%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01300"
$ws.Range("D47").Style = "Normal"
